# Weekly data refresh: a new week's price row is inserted at row 19
# (pushing the existing rows 19-23 down to 20-24), and the new row is
# populated with the latest Albahaca price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 - this shifts old rows 19:23 down to 20:24
# and extends the sheet from 23 to 24 data rows.
$ws.Rows.Item(19).Insert()

# The boilerplate columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical across
# all rows for this market/category, so copy them from the row below
# (the row that used to be row 19, now row 20) into the freshly inserted row.
$src = $ws.Rows.Item(20)
$dst = $ws.Rows.Item(19)
$src.Copy()
$dst.PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# Now overwrite row 19 with the new week's values.
$ws.Cells.Item(19, 4).Value = 44523    # D19 Fecha
$ws.Cells.Item(19, 10).Value = 250     # J19 Volumen
$ws.Cells.Item(19, 11).Value = 1400    # K19 Precio minimo
$ws.Cells.Item(19, 12).Value = 1500    # L19 Precio maximo
$ws.Cells.Item(19, 13).Value = 1450    # M19 Precio promedio ponderado
$ws.Cells.Item(19, 16).Value = 1450    # P19 Precio $/Kg
